$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = "2026-02-06 06:47:48"
$ws.Cells.Item(2, 15).Value = "-1.8 °C"
$ws.Cells.Item(3, 5).Value = "2026-02-06 06:47:51"
$ws.Cells.Item(4, 5).Value = "2026-02-06 06:47:53"
$c = $ws.Cells.Item(4, 8)
$c.NumberFormat = "@"
$c.Value = "63%"
$c.NumberFormat = "General"
$ws.Cells.Item(4, 10).Value = "993.8 hPa"
$ws.Cells.Item(4, 15).Value = "11.7 °C"
$ws.Cells.Item(5, 5).Value = "2026-02-06 06:47:56"
$ws.Cells.Item(5, 10).Value = "994.2 hPa"
$ws.Cells.Item(5, 14).Value = "5.4 °C 6:29 TU"
$ws.Cells.Item(5, 15).Value = "7.9 °C"
$ws.Cells.Item(6, 5).Value = "2026-02-06 06:47:58"
$ws.Cells.Item(6, 10).Value = "995.3 hPa"
$ws.Cells.Item(6, 14).Value = "13.4 °C 6:15 TU"
$ws.Cells.Item(6, 15).Value = "14.2 °C"
$ws.Cells.Item(7, 5).Value = "2026-02-06 06:48:01"
$c = $ws.Cells.Item(7, 8)
$c.NumberFormat = "@"
$c.Value = "68%"
$c.NumberFormat = "General"
$ws.Cells.Item(7, 10).Value = "995.1 hPa"
$ws.Cells.Item(7, 15).Value = "9.8 °C"
$ws.Cells.Item(8, 5).Value = "2026-02-06 06:48:03"
$ws.Cells.Item(8, 14).Value = "3.9 °C 6:10 TU"
$ws.Cells.Item(8, 15).Value = "5.5 °C"
$ws.Cells.Item(9, 5).Value = "2026-02-06 06:48:06"
$ws.Cells.Item(9, 14).Value = "0.0 °C 6:13 TU"
$ws.Cells.Item(9, 15).Value = "1.8 °C"
$ws.Cells.Item(10, 5).Value = "2026-02-06 06:48:08"
$ws.Cells.Item(10, 14).Value = "2.9 °C 6:09 TU"
$ws.Cells.Item(10, 15).Value = "4.8 °C"
$ws.Cells.Item(11, 5).Value = "2026-02-06 06:48:10"
$ws.Cells.Item(11, 10).Value = "996.2 hPa"
$ws.Cells.Item(11, 14).Value = "0.7 °C 6:12 TU"
$ws.Cells.Item(11, 15).Value = "3.9 °C"
$ws.Cells.Item(12, 5).Value = "2026-02-06 06:48:13"
$c = $ws.Cells.Item(12, 8)
$c.NumberFormat = "@"
$c.Value = "62%"
$c.NumberFormat = "General"
$ws.Cells.Item(12, 14).Value = "9.0 °C 6:29 TU"
$ws.Cells.Item(12, 15).Value = "11.9 °C"
$ws.Cells.Item(13, 5).Value = "2026-02-06 06:48:15"
$c = $ws.Cells.Item(13, 8)
$c.NumberFormat = "@"
$c.Value = "91%"
$c.NumberFormat = "General"
$ws.Cells.Item(13, 15).Value = "5.8 °C"
$ws.Cells.Item(14, 5).Value = "2026-02-06 06:48:18"
$ws.Cells.Item(15, 5).Value = "2026-02-06 06:48:20"
$c = $ws.Cells.Item(15, 8)
$c.NumberFormat = "@"
$c.Value = "88%"
$c.NumberFormat = "General"
$ws.Cells.Item(15, 10).Value = "994.4 hPa"
$ws.Cells.Item(15, 14).Value = "2.1 °C 6:29 TU"
$ws.Cells.Item(15, 15).Value = "5.9 °C"
$ws.Cells.Item(16, 5).Value = "2026-02-06 06:48:23"
$ws.Cells.Item(17, 5).Value = "2026-02-06 06:48:25"
$ws.Cells.Item(17, 10).Value = "997.3 hPa"
$ws.Cells.Item(17, 15).Value = "2.8 °C"
$ws.Cells.Item(18, 5).Value = "2026-02-06 06:48:28"
$ws.Cells.Item(19, 5).Value = "2026-02-06 06:48:30"
$c = $ws.Cells.Item(19, 8)
$c.NumberFormat = "@"
$c.Value = "96%"
$c.NumberFormat = "General"
$ws.Cells.Item(19, 10).Value = "997.6 hPa"
$ws.Cells.Item(19, 15).Value = "6.4 °C"
$ws.Cells.Item(20, 5).Value = "2026-02-06 06:48:32"
$ws.Cells.Item(20, 15).Value = "-2.4 °C"
$ws.Cells.Item(21, 5).Value = "2026-02-06 06:48:35"
$ws.Cells.Item(21, 10).Value = "995.3 hPa"
$ws.Cells.Item(21, 11).Value = "-0.1 MJ/m2"
$ws.Cells.Item(21, 15).Value = "4.4 °C"
$ws.Cells.Item(22, 5).Value = "2026-02-06 06:48:37"
$c = $ws.Cells.Item(22, 8)
$c.NumberFormat = "@"
$c.Value = "86%"
$c.NumberFormat = "General"
$ws.Cells.Item(22, 14).Value = "3.5 °C 6:29 TU"
$ws.Cells.Item(22, 15).Value = "7.2 °C"
$ws.Cells.Item(23, 5).Value = "2026-02-06 06:48:39"
$c = $ws.Cells.Item(23, 8)
$c.NumberFormat = "@"
$c.Value = "94%"
$c.NumberFormat = "General"
$ws.Cells.Item(23, 10).Value = "994.4 hPa"
$ws.Cells.Item(23, 14).Value = "6.1 °C 6:00 TU"
$ws.Cells.Item(23, 15).Value = "6.9 °C"
$ws.Cells.Item(24, 5).Value = "2026-02-06 06:48:42"
$ws.Cells.Item(24, 10).Value = "993.3 hPa"
$ws.Cells.Item(25, 5).Value = "2026-02-06 06:48:44"
$ws.Cells.Item(25, 10).Value = "996.5 hPa"
$ws.Cells.Item(25, 15).Value = "1.9 °C"
$ws.Cells.Item(26, 5).Value = "2026-02-06 06:48:47"
$c = $ws.Cells.Item(26, 8)
$c.NumberFormat = "@"
$c.Value = "82%"
$c.NumberFormat = "General"
$ws.Cells.Item(26, 14).Value = "-4.7 °C 6:21 TU"
$ws.Cells.Item(26, 15).Value = "-1.1 °C"
$ws.Cells.Item(27, 5).Value = "2026-02-06 06:48:49"
$c = $ws.Cells.Item(27, 8)
$c.NumberFormat = "@"
$c.Value = "98%"
$c.NumberFormat = "General"
$ws.Cells.Item(27, 10).Value = "994.1 hPa"
$ws.Cells.Item(27, 15).Value = "7.1 °C"
$ws.Cells.Item(28, 5).Value = "2026-02-06 06:48:52"
$c = $ws.Cells.Item(28, 8)
$c.NumberFormat = "@"
$c.Value = "91%"
$c.NumberFormat = "General"
$ws.Cells.Item(28, 10).Value = "997.5 hPa"
$ws.Cells.Item(28, 14).Value = "-0.9 °C 6:21 TU"
$ws.Cells.Item(28, 15).Value = "2.1 °C"
$ws.Cells.Item(29, 5).Value = "2026-02-06 06:48:54"
$c = $ws.Cells.Item(29, 8)
$c.NumberFormat = "@"
$c.Value = "68%"
$c.NumberFormat = "General"
$ws.Cells.Item(29, 14).Value = "5.7 °C 6:25 TU"
$ws.Cells.Item(29, 15).Value = "10.4 °C"
$ws.Cells.Item(30, 5).Value = "2026-02-06 06:48:56"
$c = $ws.Cells.Item(30, 8)
$c.NumberFormat = "@"
$c.Value = "72%"
$c.NumberFormat = "General"
$ws.Cells.Item(30, 15).Value = "-3.7 °C"
$ws.Cells.Item(31, 5).Value = "2026-02-06 06:48:59"
$ws.Cells.Item(31, 10).Value = "997.2 hPa"
$ws.Cells.Item(32, 5).Value = "2026-02-06 06:49:01"
$ws.Cells.Item(32, 10).Value = "995.7 hPa"
$ws.Cells.Item(32, 15).Value = "14.5 °C"
$ws.Cells.Item(33, 5).Value = "2026-02-06 06:49:04"
$ws.Cells.Item(33, 14).Value = "4.3 °C 6:03 TU"
$ws.Cells.Item(33, 15).Value = "6.3 °C"
$ws.Cells.Item(34, 5).Value = "2026-02-06 06:49:06"
$c = $ws.Cells.Item(34, 8)
$c.NumberFormat = "@"
$c.Value = "83%"
$c.NumberFormat = "General"
$ws.Cells.Item(34, 14).Value = "2.7 °C 6:29 TU"
$ws.Cells.Item(34, 15).Value = "6.9 °C"
$ws.Cells.Item(35, 5).Value = "2026-02-06 06:49:09"
$ws.Cells.Item(35, 14).Value = "-3.4 °C 6:14 TU"
$ws.Cells.Item(35, 15).Value = "-3.2 °C"
$ws.Cells.Item(36, 5).Value = "2026-02-06 06:49:11"
$c = $ws.Cells.Item(36, 8)
$c.NumberFormat = "@"
$c.Value = "69%"
$c.NumberFormat = "General"
$ws.Cells.Item(36, 10).Value = "997.1 hPa"
$ws.Cells.Item(36, 14).Value = "7.9 °C 6:19 TU"
$ws.Cells.Item(36, 15).Value = "11.0 °C"
